$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Swap the match data (columns F:V) between row 29 and row 31.
#    Columns A:E (index, pais, torneio, temporada, data_partida)
#    stay exactly as they were.
# ---------------------------------------------------------------
$swapCols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$row29vals = @{}
$row31vals = @{}
foreach ($c in $swapCols) {
    $row29vals[$c] = $ws.Range($c + "29").Value2
    $row31vals[$c] = $ws.Range($c + "31").Value2
}
foreach ($c in $swapCols) {
    $ws.Range($c + "29").Value2 = $row31vals[$c]
    $ws.Range($c + "31").Value2 = $row29vals[$c]
}

# ---------------------------------------------------------------
# 2) Append 5 new match rows (65-69) after the current last row
#    (row 64), copying formatting from row 64 first.
# ---------------------------------------------------------------
for ($i = 65; $i -le 69; $i++) {
    $ws.Range("A64:V64").Copy($ws.Range("A" + $i + ":V" + $i))
}

$newRows = @(
    @{ Row=65; A=64; E=45192.66666666666; F="Avia Swidnik";     G=6; H="Czarni Polaniec";    I=2;
       J=1.37; K="22/09/2023 03:13"; L=1.39; M="23/09/2023 15:44";
       N=4.35; O="22/09/2023 03:13"; P=4.49; Q="23/09/2023 15:48";
       R=5.22; S="22/09/2023 03:13"; T=6.11; U="23/09/2023 15:44";
       V="https://www.betexplorer.com/football/poland/iii-liga-group-iv/avia-swidnik-czarni-polaniec/Mw3Z2svq/" },

    @{ Row=66; A=65; E=45192.66666666666; F="Chelmianka Chelm"; G=2; H="Star Starachowice";  I=0;
       J=2.1;  K="22/09/2023 03:13"; L=2.6;  M="23/09/2023 15:08";
       N=3.2;  O="22/09/2023 03:13"; P=3.37; Q="23/09/2023 15:08";
       R=2.8;  S="22/09/2023 03:13"; T=2.36; U="23/09/2023 15:08";
       V="https://www.betexplorer.com/football/poland/iii-liga-group-iv/chelmianka-chelm-star-starachowice/KUG77NVR/" },

    @{ Row=67; A=66; E=45192.66666666666; F="Orleta Radzyn";    G=0; H="Podhale Nowy Targ";  I=0;
       J=3.04; K="22/09/2023 03:13"; L=2.77; M="23/09/2023 15:13";
       N=3.4;  O="22/09/2023 03:13"; P=3.6;  Q="23/09/2023 15:13";
       R=1.9;  S="22/09/2023 03:13"; T=2.14; U="23/09/2023 15:13";
       V="https://www.betexplorer.com/football/poland/iii-liga-group-iv/orleta-radzyn-podhale-nowy-targ/CSen0uO1/" },

    @{ Row=68; A=67; E=45192.66666666666; F="Unia Tarnow";      G=1; H="Wisloka Debica";     I=1;
       J=2.14; K="22/09/2023 03:13"; L=2.25; M="23/09/2023 15:51";
       N=3.22; O="22/09/2023 03:13"; P=3.49; Q="23/09/2023 15:51";
       R=2.72; S="22/09/2023 03:13"; T=2.67; U="23/09/2023 15:44";
       V="https://www.betexplorer.com/football/poland/iii-liga-group-iv/unia-tarnow-wisloka-debica/hvmP7is8/" },

    @{ Row=69; A=68; E=45193.79166666666; F="Wiazownica";       G=1; H="Siarka Tarnobrzeg";  I=1;
       J=4.14; K="24/09/2023 13:42"; L=3.24; M="24/09/2023 18:44";
       N=3.85; O="24/09/2023 13:42"; P=3.52; Q="24/09/2023 18:44";
       R=1.61; S="24/09/2023 13:42"; T=1.94; U="24/09/2023 18:44";
       V="https://www.betexplorer.com/football/poland/iii-liga-group-iv/wiazownica-siarka-tarnobrzeg/42hv21gk/" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = $r.A
    $ws.Cells.Item($row, 2).Value2  = "poland"
    $ws.Cells.Item($row, 3).Value2  = "iii-liga-group-iv"
    $ws.Cells.Item($row, 4).Value2  = "2023-2024"
    $ws.Cells.Item($row, 5).Value2  = $r.E
    $ws.Cells.Item($row, 6).Value2  = $r.F
    $ws.Cells.Item($row, 7).Value2  = $r.G
    $ws.Cells.Item($row, 8).Value2  = $r.H
    $ws.Cells.Item($row, 9).Value2  = $r.I
    $ws.Cells.Item($row, 10).Value2 = $r.J
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $r.N
    $ws.Cells.Item($row, 15).Value2 = $r.O
    $ws.Cells.Item($row, 16).Value2 = $r.P
    $ws.Cells.Item($row, 17).Value2 = $r.Q
    $ws.Cells.Item($row, 18).Value2 = $r.R
    $ws.Cells.Item($row, 19).Value2 = $r.S
    $ws.Cells.Item($row, 20).Value2 = $r.T
    $ws.Cells.Item($row, 21).Value2 = $r.U
    $ws.Cells.Item($row, 22).Value2 = $r.V
}
